$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2535284767066628
$ws.Range("D2").Value = 0.03384184930752099
$ws.Range("E2").Value = 0.1460973905544876
$ws.Range("F2").Value = 0.570786126120403
$ws.Range("G2").Value = 0.4119813430658041
$ws.Range("H2").Value = 0.5818873947717478
$ws.Range("I2").Value = 0.6100275392125774
$ws.Range("K2").Value = 1.8271088237085
$ws.Range("L2").Value = 0.1657999108884667
$ws.Range("O2").Value = 1.935801470792597
$ws.Range("C3").Value = 0.2459551058914968
$ws.Range("D3").Value = 0.03077883211933852
$ws.Range("E3").Value = 0.1420438851354504
$ws.Range("F3").Value = 0.5762823226489218
$ws.Range("G3").Value = 0.4187710497625545
$ws.Range("H3").Value = 0.5898434115779878
$ws.Range("I3").Value = 0.615983827839262
$ws.Range("K3").Value = 1.604884620796497
$ws.Range("L3").Value = 0.1614548716992559
$ws.Range("O3").Value = 1.966560233321019
$ws.Range("C4").Value = 0.2413992630319797
$ws.Range("D4").Value = 0.02888790721193146
$ws.Range("E4").Value = 0.139626942099639
$ws.Range("F4").Value = 0.5801803780047905
$ws.Range("G4").Value = 0.4234010830165076
$ws.Range("H4").Value = 0.5950985632668804
$ws.Range("I4").Value = 0.6201408263452244
$ws.Range("K4").Value = 1.467890819799038
$ws.Range("L4").Value = 0.1588840840587693
$ws.Range("O4").Value = 1.987188871812862
$ws.Range("C5").Value = 0.2395665547145143
$ws.Range("D5").Value = 0.02811482424070277
$ws.Range("E5").Value = 0.1386601477558322
$ws.Range("F5").Value = 0.5819001913013295
$ws.Range("G5").Value = 0.4254034280386705
$ws.Range("H5").Value = 0.5973330628598248
$ws.Range("I5").Value = 0.6219603016009323
$ws.Range("K5").Value = 1.411931352451404
$ws.Range("L5").Value = 0.1578608908672763
$ws.Range("O5").Value = 1.996032401628042
$ws.Range("C6").Value = 0.2392636785209419
$ws.Range("D6").Value = 0.02798630406406488
$ws.Range("E6").Value = 0.1385007085542362
$ws.Range("F6").Value = 0.5821936883713974
$ws.Range("G6").Value = 0.4257428846658158
$ws.Range("H6").Value = 0.5977097131395439
$ws.Range("I6").Value = 0.6222699954419788
$ws.Range("K6").Value = 1.402631394891841
$ws.Range("L6").Value = 0.1576924658816878
$ws.Range("O6").Value = 1.997527238743686
$ws.Range("C7").Value = 0.2413744497985988
$ws.Range("D7").Value = 0.02887749125844863
$ws.Range("E7").Value = 0.1396138301092265
$ws.Range("F7").Value = 0.5802030406566132
$ws.Range("G7").Value = 0.4234276199195861
$ws.Range("H7").Value = 0.5951283221456336
$ws.Range("I7").Value = 0.6201648567131883
$ws.Range("K7").Value = 1.467136665271994
$ws.Range("L7").Value = 0.1588701860117609
$ws.Range("O7").Value = 1.987306369915842
$ws.Range("C8").Value = 0.2508976991587843
$ws.Range("D8").Value = 0.03278787358417645
$ws.Range("E8").Value = 0.1446848380197423
$ws.Range("F8").Value = 0.5725724204389095
$ws.Range("G8").Value = 0.4142264965257993
$ws.Range("H8").Value = 0.5845537574252049
$ws.Range("I8").Value = 0.6119774004569223
$ws.Range("K8").Value = 1.750601954362139
$ws.Range("L8").Value = 0.1642815964343853
$ws.Range("O8").Value = 1.946044805754156
$ws.Range("C9").Value = 0.270315608609593
$ws.Range("D9").Value = 0.04037301997783516
$ws.Range("E9").Value = 0.1551985438235519
$ws.Range("F9").Value = 0.5617744437591057
$ws.Range("G9").Value = 0.3998591919051862
$ws.Range("H9").Value = 0.5667576644274845
$ws.Range("I9").Value = 0.5998976657291237
$ws.Range("K9").Value = 2.301976376583752
$ws.Range("L9").Value = 0.1756639644864748
$ws.Range("O9").Value = 1.879002855878355
$ws.Range("C10").Value = 0.2850299678272279
$ws.Range("D10").Value = 0.04589300757150738
$ws.Range("E10").Value = 0.1632694884550219
$ws.Range("F10").Value = 0.5563989767970554
$ws.Range("G10").Value = 0.3915682528090869
$ws.Range("H10").Value = 0.555480588394154
$ws.Range("I10").Value = 0.5934607186227723
$ws.Range("K10").Value = 2.704156078571771
$ws.Range("L10").Value = 0.1844979173025223
$ws.Range("O10").Value = 1.838264972841259
$ws.Range("C11").Value = 0.291820284780357
$ws.Range("D11").Value = 0.04839230633238856
$ws.Range("E11").Value = 0.1670163435748435
$ws.Range("F11").Value = 0.554512878323429
$ws.Range("G11").Value = 0.3882933096932888
$ws.Range("H11").Value = 0.5507417886598631
$ws.Range("I11").Value = 0.5910648535649656
$ws.Range("K11").Value = 2.886451817507805
$ws.Range("L11").Value = 0.1886194561281798
$ws.Range("O11").Value = 1.821594885404025
$ws.Range("C12").Value = 0.2944053908640285
$ws.Range("D12").Value = 0.04933698696123656
$ws.Range("E12").Value = 0.1684459862294361
$ws.Range("F12").Value = 0.5538793881676654
$ws.Range("G12").Value = 0.3871250033543348
$ws.Range("H12").Value = 0.5490036800710456
$ws.Range("I12").Value = 0.5902343939948764
$ws.Range("K12").Value = 2.955384288256255
$ws.Range("L12").Value = 0.1901949844811526
$ws.Range("O12").Value = 1.815551163449257
$ws.Range("C13").Value = 0.2938480322737291
$ws.Range("D13").Value = 0.04913361179769993
$ws.Range("E13").Value = 0.168137607921075
$ws.Range("F13").Value = 0.554012225884037
$ws.Range("G13").Value = 0.3873734170014771
$ws.Range("H13").Value = 0.5493755038157602
$ws.Range("I13").Value = 0.5904098284904222
$ws.Range("K13").Value = 2.940542918143308
$ws.Range("L13").Value = 0.1898550085870028
$ws.Range("O13").Value = 1.816840810309657
$ws.Range("C14").Value = 0.2920326880465041
$ws.Range("D14").Value = 0.04847006116334285
$ws.Range("E14").Value = 0.1671337451425501
$ws.Range("F14").Value = 0.5544591408561033
$ws.Range("G14").Value = 0.3881957501415556
$ws.Range("H14").Value = 0.5505976629656217
$ws.Range("I14").Value = 0.5909949904731988
$ws.Range("K14").Value = 2.892124941739837
$ws.Range("L14").Value = 0.188748779417665
$ws.Range("O14").Value = 1.821092269432313
$ws.Range("C15").Value = 0.2909225250924976
$ws.Range("D15").Value = 0.04806338771639673
$ws.Range("E15").Value = 0.1665202540642596
$ws.Range("F15").Value = 0.5547434125818
$ws.Range("G15").Value = 0.3887088217013144
$ws.Range("H15").Value = 0.5513536155348291
$ws.Range("I15").Value = 0.5913634285349971
$ws.Range("K15").Value = 2.862454515426521
$ws.Range("L15").Value = 0.1880731080396174
$ws.Range("O15").Value = 1.823731461039898
$ws.Range("C16").Value = 0.2845881379471109
$ws.Range("D16").Value = 0.04572943077650393
$ws.Range("E16").Value = 0.1630261349589119
$ws.Range("F16").Value = 0.5565335181142714
$ws.Range("G16").Value = 0.3917923083401931
$ws.Range("H16").Value = 0.5557981621261092
$ws.Range("I16").Value = 0.5936280276759121
$ws.Range("K16").Value = 2.692229011103336
$ws.Range("L16").Value = 0.184230635391728
$ws.Range("O16").Value = 1.839391965032291
$ws.Range("C17").Value = 0.2807268575838009
$ws.Range("D17").Value = 0.04429456824213673
$ws.Range("E17").Value = 0.1609018715115127
$ws.Range("F17").Value = 0.5577751487327731
$ws.Range("G17").Value = 0.3938114346983852
$ws.Range("H17").Value = 0.5586250270463324
$ws.Range("I17").Value = 0.5951538104156882
$ws.Range("K17").Value = 2.587629564765336
$ws.Range("L17").Value = 0.1818997584452973
$ws.Range("O17").Value = 1.849476848222437
$ws.Range("C18").Value = 0.2785150597997301
$ws.Range("D18").Value = 0.04346816865958658
$ws.Range("E18").Value = 0.1596871455363029
$ws.Range("F18").Value = 0.5585419201853981
$ws.Range("G18").Value = 0.3950195041585687
$ws.Range("H18").Value = 0.56028777941453
$ws.Range("I18").Value = 0.5960814910821028
$ws.Range("K18").Value = 2.527405070591783
$ws.Range("L18").Value = 0.1805687884896088
$ws.Range("O18").Value = 1.855452577018767
$ws.Range("C19").Value = 0.27776775181772
$ws.Range("D19").Value = 0.04318817598618807
$ws.Range("E19").Value = 0.1592770799698613
$ws.Range("F19").Value = 0.5588105633505194
$ws.Range("G19").Value = 0.3954365471375141
$ws.Range("H19").Value = 0.5608570781070057
$ws.Range("I19").Value = 0.5964041832081648
$ws.Range("K19").Value = 2.507003635745718
$ws.Range("L19").Value = 0.1801198092985317
$ws.Range("O19").Value = 1.857505906525631
$ws.Range("C20").Value = 0.2811369558041577
$ws.Range("D20").Value = 0.04444742654501965
$ws.Range("E20").Value = 0.1611272690818879
$ws.Range("F20").Value = 0.5576375269189526
$ws.Range("G20").Value = 0.3935916566668709
$ws.Range("H20").Value = 0.5583202914357912
$ws.Range("I20").Value = 0.5949862021872576
$ws.Range("K20").Value = 2.598770767733868
$ws.Range("L20").Value = 0.1821468815924163
$ws.Range("O20").Value = 1.8483851572784
$ws.Range("C21").Value = 0.2925655263266265
$ws.Range("D21").Value = 0.04866501000303458
$ws.Range("E21").Value = 0.1674283114903758
$ws.Range("F21").Value = 0.5543256774640497
$ws.Range("G21").Value = 0.3879522577645105
$ws.Range("H21").Value = 0.5502371543375375
$ws.Range("I21").Value = 0.5908210278038268
$ws.Range("K21").Value = 2.906349195698965
$ws.Range("L21").Value = 0.1890733043965156
$ws.Range("O21").Value = 1.819836206357849
$ws.Range("C22").Value = 0.300114898542688
$ws.Range("D22").Value = 0.05141122076010163
$ws.Range("E22").Value = 0.1716092785306387
$ws.Range("F22").Value = 0.5526319201957364
$ws.Range("G22").Value = 0.3846855400554645
$ws.Range("H22").Value = 0.5452829908304651
$ws.Range("I22").Value = 0.5885466259432377
$ws.Range("K22").Value = 3.10679061450378
$ws.Range("L22").Value = 0.1936863308723673
$ws.Range("O22").Value = 1.80274557843039
$ws.Range("C23").Value = 0.29607837045927
$ws.Range("D23").Value = 0.04994647052868117
$ws.Range("E23").Value = 0.1693720795497953
$ws.Range("F23").Value = 0.5534927342310709
$ws.Range("G23").Value = 0.3863905740122391
$ws.Range("H23").Value = 0.5478970121798525
$ws.Range("I23").Value = 0.5897194617646591
$ws.Range("K23").Value = 2.99986565600625
$ws.Range("L23").Value = 0.1912163870816954
$ws.Range("O23").Value = 1.811723337113818
$ws.Range("C24").Value = 0.2809515251483106
$ws.Range("D24").Value = 0.04437832390578933
$ws.Range("E24").Value = 0.1610253464621749
$ws.Range("F24").Value = 0.5576995808557044
$ws.Range("G24").Value = 0.3936908711759912
$ws.Range("H24").Value = 0.5584579454469463
$ws.Range("I24").Value = 0.595061820619641
$ws.Range("K24").Value = 2.593734106391764
$ws.Range("L24").Value = 0.1820351289277369
$ws.Range("O24").Value = 1.848878157017552
$ws.Range("C25").Value = 0.2649835254747188
$ws.Range("D25").Value = 0.03833015877427215
$ws.Range("E25").Value = 0.1522934084138896
$ws.Range("F25").Value = 0.5642477198940483
$ws.Range("G25").Value = 0.4033499237718416
$ws.Range("H25").Value = 0.5712566452279901
$ws.Range("I25").Value = 0.602738433687307
$ws.Range("K25").Value = 2.153314696998848
$ws.Range("L25").Value = 0.1880731080396174
$ws.Range("O25").Value = 1.895648059678535
